# Updates cryptos list: new Price (col D) and Volume(1h) (col E) values.
# Values are prefixed with a leading apostrophe so Excel stores them as
# literal text (matching the original inlineStr cells) instead of coercing
# numeric-looking strings (e.g. "228.70") into floating point numbers. The
# Style reset afterwards clears the quotePrefix formatting flag that the
# apostrophe entry leaves behind, so only the cell VALUE changes - matching
# the target diff which shows no style/attribute changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.331.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.17%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.045.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.44%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'228.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.82%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  -1.47%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'56.93"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.81%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.385"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.72%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.0789"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.52%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -1.85%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'14.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.08%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'2.344.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.51%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'20.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.08%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.756"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.38%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'5.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.78%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.057.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.39%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'37.196.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.39%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'6.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.51%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'69.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.67%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.0₃0827"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.45%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'225.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.44%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.09%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.15%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -4.49%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'9.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.14%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'168.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.17%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.128"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -6.61%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'18.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.50%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.44%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  -2.16%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'4.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.94%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.0613"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.87%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'4.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.18%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.79%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +1.53%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.02%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -4.51%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -2.63%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -4.87%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'17.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.09%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -1.16%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.472.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.72%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.0942"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.19%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'96.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.77%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.70%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -4.03%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'3.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.77%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'7.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.78%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -2.12%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'2.229.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.65%  "
$ws.Range("E51").Style = "Normal"
